# "classification on the file level" -- the sheet-level "classification" and
# "position" columns are retired from the sessions header row; "classification"
# and "position" are now tracked per-file (fclassification_1 / fclassification_2,
# fposition_1 / fposition_2). Two header cells ("fposition_1clip_out_1" and
# "fname_2fposition_2") were also mangled/merged and are split back into their
# proper individual columns. Net effect: the 27-column header (A:AA) becomes a
# 29-column header (A:AC).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sessions")

$headers = @(
  "name","key","date","participantID","top","pilot","exclusion","setting",
  "country","state","language","consent","condition","group","tasks",
  "transcode_options","filepath","file_1","fname_1","fposition_1",
  "fclassification_1","clip_out_1","clip_in_1","file_2","fname_2",
  "fposition_2","fclassification_2","clip_out_2","clip_in_2"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Rebuild the dropdown (list) data validations against their new column
# positions now that the header has been reshuffled.
$ws.Cells.Validation.Delete()

# exclusion (col G) -- unchanged
$ws.Range("G2:G1001").Validation.Add(3, 1, 1, '"Did not meet inclusion criteria,Procedural/experimenter error,Withdrew/fussy/tired,Outlier"')

# fclassification_1 (col U) -- new, file-level classification
$ws.Range("U2:U1001").Validation.Add(3, 1, 1, '"None,PRIVATE,SHARED,EXCERPTS,PUBLIC"')

# fclassification_2 (col AA) -- new, file-level classification
$ws.Range("AA2:AA1001").Validation.Add(3, 1, 1, '"None,PRIVATE,SHARED,EXCERPTS,PUBLIC"')

# setting (col H, was I)
$ws.Range("H2:H1001").Validation.Add(3, 1, 1, '"Lab,Home,Classroom,Outdoor,Clinic"')

# state (col J, was K)
$ws.Range("J2:J1001").Validation.Add(3, 1, 1, '"AL,AK,AZ,AR,CA,CO,CT,DE,DC,FL,GA,HI,ID,IL,IN,IA,KS,KY,LA,ME,MT,NE,NV,NH,NJ,NM,NY,NC,ND,OH,OK,OR,MD,MA,MI,MN,MS,MO,PA,RI,SC,SD,TN,TX,UT,VT,VA,WA,WV,WI,WY"')

# consent (col L, was M)
$ws.Range("L2:L1001").Validation.Add(3, 1, 1, '"None,PRIVATE,SHARED,EXCERPTS,PUBLIC"')
